# Day 7 night written - add Raiders arriving sequence and related conversations to Lines sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lines")

# --- Write new row data (rows 2150-2303) ---
$ws.Range("A2150").Value = "CONVERSATION"
$ws.Range("B2150").Value = "Raiders arriving"
$ws.Range("D2150").Value = 7
$ws.Range("A2151").Value = "Raiders"
$ws.Range("B2151").Value = "HEY! WE KNOW YOU'RE IN THERE! OPEN UP!"
$ws.Range("D2151").Value = 7
$ws.Range("A2152").Value = "Raiders"
$ws.Range("B2152").Value = "OPEN UP OR WE ARE BUSTING OUR WAY IN!"
$ws.Range("D2152").Value = 7
$ws.Range("A2153").Value = "IF"
$ws.Range("B2153").Value = "Bob"
$ws.Range("A2154").Value = "Bob"
$ws.Range("B2154").Value = "Oh shit."
$ws.Range("D2154").Value = 7
$ws.Range("A2155").Value = "ENDIF"
$ws.Range("A2156").Value = "Raiders"
$ws.Range("B2156").Value = "WE WON'T HURT YOU! OPEN THE GODDAMN DOOR."
$ws.Range("D2156").Value = 7
$ws.Range("A2157").Value = "END"
$ws.Range("A2160").Value = "CONVERSATION"
$ws.Range("B2160").Value = "Look through peephole"
$ws.Range("D2160").Value = 7
$ws.Range("A2161").Value = "Player"
$ws.Range("B2161").Value = "That is a lot of people. At least five. They have guns."
$ws.Range("D2161").Value = 7
$ws.Range("A2162").Value = "END"
$ws.Range("A2165").Value = "CONVERSATION"
$ws.Range("B2165").Value = "Stay silent"
$ws.Range("D2165").Value = 7
$ws.Range("A2166").Value = "Raiders"
$ws.Range("B2166").Value = "WE CAN SEE YOU PUSSIES IN THERE."
$ws.Range("D2166").Value = 7
$ws.Range("A2167").Value = "END"
$ws.Range("A2170").Value = "CONVERSATION"
$ws.Range("B2170").Value = "Who are you"
$ws.Range("D2170").Value = 7
$ws.Range("A2171").Value = "Raiders"
$ws.Range("B2171").Value = "OPEN THE FUCKING DOOR AND FIND OUT."
$ws.Range("D2171").Value = 7
$ws.Range("A2172").Value = "END"
$ws.Range("A2175").Value = "CONVERSATION"
$ws.Range("B2175").Value = "What do you want"
$ws.Range("D2175").Value = 7
$ws.Range("A2176").Value = "Raiders"
$ws.Range("B2176").Value = "I WANT TO KILL YOU FOR ASKING SO MANY QUESTIONS."
$ws.Range("D2176").Value = 7
$ws.Range("A2177").Value = "END"
$ws.Range("A2180").Value = "CONVERSATION"
$ws.Range("B2180").Value = "Why should I let you in"
$ws.Range("D2180").Value = 7
$ws.Range("A2181").Value = "Raiders"
$ws.Range("B2181").Value = "WE WILL COME IN WHETHER YOU OPEN THIS DOOR OR NOT. IF YOU LET US IN, WE MIGHT NOT KILL YOU ALL."
$ws.Range("D2181").Value = 7
$ws.Range("A2182").Value = "END"
$ws.Range("A2185").Value = "CONVERSATION"
$ws.Range("B2185").Value = "How can you be helpful"
$ws.Range("D2185").Value = 7
$ws.Range("A2186").Value = "Raiders"
$ws.Range("B2186").Value = "*shoots bullets through door*"
$ws.Range("D2186").Value = 7
$ws.Range("A2187").Value = "IF"
$ws.Range("B2187").Value = "Jessica"
$ws.Range("A2188").Value = "Jessica"
$ws.Range("B2188").Value = "AH!"
$ws.Range("D2188").Value = 7
$ws.Range("A2189").Value = "ENDIF"
$ws.Range("A2190").Value = "IF"
$ws.Range("B2190").Value = "Bob"
$ws.Range("A2191").Value = "Bob"
$ws.Range("B2191").Value = "Holy shit."
$ws.Range("D2191").Value = 7
$ws.Range("A2192").Value = "ENDIF"
$ws.Range("A2193").Value = "END"
$ws.Range("A2196").Value = "CONVERSATION"
$ws.Range("B2196").Value = "Bobs opinion of raiders"
$ws.Range("A2197").Value = "Bob"
$ws.Range("B2197").Value = "Barricade the fucking door. If you let them in, we are dead."
$ws.Range("A2198").Value = "END"
$ws.Range("A2200").Value = "CONVERSATION"
$ws.Range("B2200").Value = "Jessicas opinion of raiders"
$ws.Range("A2201").Value = "Jessica"
$ws.Range("B2201").Value = "Um.. oh god. Maybe you should let them in? Maybe they will just leave..."
$ws.Range("A2202").Value = "IF"
$ws.Range("B2202").Value = "BOB"
$ws.Range("A2203").Value = "Bob"
$ws.Range("B2203").Value = "Shut the fuck up. Do not let them in."
$ws.Range("A2204").Value = "ENDIF"
$ws.Range("A2205").Value = "END"
$ws.Range("A2208").Value = "CONVERSATION"
$ws.Range("B2208").Value = "Violets opinion of raiders"
$ws.Range("A2209").Value = "Violet"
$ws.Range("B2209").Value = "(curling fists) ..."
$ws.Range("A2210").Value = "END"
$ws.Range("A2213").Value = "CONVERSATION"
$ws.Range("B2213").Value = "Hals opinion of raiders"
$ws.Range("A2214").Value = "Hal"
$ws.Range("B2214").Value = "These fuckers better get ready. We aren't going down without a fight."
$ws.Range("A2215").Value = "END"
$ws.Range("A2218").Value = "CONVERSATION"
$ws.Range("B2218").Value = "Sals opinion of raiders"
$ws.Range("A2219").Value = "Sal"
$ws.Range("B2219").Value = "*whispering* I hope your door is strong."
$ws.Range("A2220").Value = "END"
$ws.Range("A2224").Value = "CONVERSATION"
$ws.Range("B2224").Value = "Keep raiders out"
$ws.Range("A2225").Value = "Player"
$ws.Range("B2225").Value = "Get the fuck out of here."
$ws.Range("A2226").Value = "Raiders"
$ws.Range("B2226").Value = "YOU ASKED FOR IT."
$ws.Range("A2227").Value = "Door"
$ws.Range("B2227").Value = "*slamming sounds*"
$ws.Range("A2228").Value = "Raiders"
$ws.Range("B2228").Value = "(muffled) Back up guys..."
$ws.Range("A2229").Value = "None"
$ws.Range("B2229").Value = "Shots ring out."
$ws.Range("A2230").Value = "None"
$ws.Range("B2230").Value = "You are hit."
$ws.Range("A2231").Value = "IF"
$ws.Range("B2231").Value = "Hal"
$ws.Range("A2232").Value = "Hal"
$ws.Range("B2232").Value = "Shit, you're bleeding."
$ws.Range("A2233").Value = "ELIF"
$ws.Range("B2233").Value = "Bob"
$ws.Range("A2234").Value = "Bob"
$ws.Range("B2234").Value = "Put some pressure on that."
$ws.Range("A2235").Value = "ENDIF"
$ws.Range("A2236").Value = "Raiders"
$ws.Range("B2236").Value = "HOW DO YOU FUCKING LIKE THAT?"
$ws.Range("A2237").Value = "None"
$ws.Range("B2237").Value = "More shots ring out, hitting nothing."
$ws.Range("A2238").Value = "Raiders"
$ws.Range("B2238").Value = "WE WILL BE BACK YOU ASSHOLE. YOU'RE DEAD!"
$ws.Range("A2239").Value = "END"
$ws.Range("A2243").Value = "CONVERSATION"
$ws.Range("B2243").Value = "Let raiders in"
$ws.Range("A2244").Value = "Player"
$ws.Range("B2244").Value = "Okay, I'll open the door."
$ws.Range("A2245").Value = "IF"
$ws.Range("B2245").Value = "Jessica"
$ws.Range("A2246").Value = "Jessica"
$ws.Range("B2246").Value = "DON'T-"
$ws.Range("A2247").Value = "ENDIF"
$ws.Range("A2248").Value = "IF"
$ws.Range("B2248").Value = "Bob"
$ws.Range("A2249").Value = "Bob"
$ws.Range("B2249").Value = "HEY!"
$ws.Range("A2250").Value = "ENDIF"
$ws.Range("A2251").Value = "IF"
$ws.Range("B2251").Value = "Violet"
$ws.Range("A2252").Value = "Violet"
$ws.Range("B2252").Value = "*yelling*"
$ws.Range("A2253").Value = "ENDIF"
$ws.Range("A2254").Value = "IF"
$ws.Range("B2254").Value = "Hal"
$ws.Range("A2255").Value = "Hal"
$ws.Range("B2255").Value = "STOP!"
$ws.Range("A2256").Value = "ENDIF"
$ws.Range("A2257").Value = "None"
$ws.Range("B2257").Value = "The raiders come inside."
$ws.Range("A2258").Value = "Raiders"
$ws.Range("B2258").Value = "Good. You did the right thing, huh?"
$ws.Range("A2259").Value = "IF"
$ws.Range("B2259").Value = "Jessica"
$ws.Range("A2260").Value = "None"
$ws.Range("B2260").Value = "Yells echo as the raiders drag everyone out."
$ws.Range("A2261").Value = "Raiders"
$ws.Range("B2261").Value = "We'll get rid of all your little problems. Little revolutionary scum, all of them."
$ws.Range("A2262").Value = "ELIF"
$ws.Range("B2262").Value = "Bob"
$ws.Range("A2263").Value = "None"
$ws.Range("B2263").Value = "Yells echo as the raiders drag everyone out."
$ws.Range("A2264").Value = "Raiders"
$ws.Range("B2264").Value = "We'll get rid of all your little problems. Little revolutionary scum, all of them."
$ws.Range("A2265").Value = "ELIF"
$ws.Range("B2265").Value = "Violet"
$ws.Range("A2266").Value = "None"
$ws.Range("B2266").Value = "Yells echo as the raiders drag everyone out."
$ws.Range("A2267").Value = "Raiders"
$ws.Range("B2267").Value = "We'll get rid of all your little problems. Little revolutionary scum, all of them."
$ws.Range("A2268").Value = "ELIF"
$ws.Range("B2268").Value = "Hal"
$ws.Range("A2269").Value = "None"
$ws.Range("B2269").Value = "Yells echo as the raiders drag everyone out."
$ws.Range("A2270").Value = "Raiders"
$ws.Range("B2270").Value = "We'll get rid of all your little problems. Little revolutionary scum, all of them."
$ws.Range("A2271").Value = "ENDIF"
$ws.Range("A2272").Value = "Raiders"
$ws.Range("B2272").Value = "Glory to Brasnia."
$ws.Range("A2273").Value = "None"
$ws.Range("B2273").Value = "The raiders leave and close the door behind them."
$ws.Range("A2274").Value = "//"
$ws.Range("B2274").Value = "I love not implementing Ors"
$ws.Range("A2275").Value = "IF"
$ws.Range("B2275").Value = "Jessica"
$ws.Range("A2276").Value = "GOTO"
$ws.Range("B2276").Value = "Execute party members"
$ws.Range("A2277").Value = "ELIF"
$ws.Range("B2277").Value = "Bob"
$ws.Range("A2278").Value = "GOTO"
$ws.Range("B2278").Value = "Execute party members"
$ws.Range("A2279").Value = "ELIF"
$ws.Range("B2279").Value = "Violet"
$ws.Range("A2280").Value = "GOTO"
$ws.Range("B2280").Value = "Execute party members"
$ws.Range("A2281").Value = "ELIF"
$ws.Range("B2281").Value = "Hal"
$ws.Range("A2282").Value = "GOTO"
$ws.Range("B2282").Value = "Execute party members"
$ws.Range("A2283").Value = "ENDIF"
$ws.Range("A2284").Value = "GOTO"
$ws.Range("B2284").Value = "Day7Over"
$ws.Range("A2285").Value = "END"
$ws.Range("A2291").Value = "CONVERSATION"
$ws.Range("B2291").Value = "Execute party members"
$ws.Range("A2292").Value = "None"
$ws.Range("B2292").Value = "You hear pleading and crying."
$ws.Range("A2293").Value = "None"
$ws.Range("B2293").Value = "Shots ring out."
$ws.Range("A2294").Value = "IF"
$ws.Range("B2294").Value = "HasCar"
$ws.Range("A2295").Value = "None"
$ws.Range("B2295").Value = "You hear nothing but the car starting up and driving away."
$ws.Range("A2296").Value = "ELSE"
$ws.Range("A2297").Value = "None"
$ws.Range("B2297").Value = "You hear nothing but the wind."
$ws.Range("A2298").Value = "ENDIF"
$ws.Range("A2299").Value = "GOTO"
$ws.Range("B2299").Value = "Day7Over"
$ws.Range("A2300").Value = "END"
$ws.Range("A2302").Value = "CONVERSATION"
$ws.Range("B2302").Value = "Day7Over"
$ws.Range("A2303").Value = "END"

# --- Rebuild the conditional formatting range (sqref) for the main rule block ---
# so it matches the fragmented per-block ranges Excel produces once the new
# rows (with differing used-column extents) are inserted under the existing rule.
$cfRanges = @(
    "A1957:H1959",
    "A1956",
    "C1956:H1956",
    "A1961:H1964",
    "A1960",
    "C1960:H1960",
    "A1966:H1968",
    "A1965",
    "C1965:H1965",
    "A1969",
    "C1969:H1969",
    "A1970:H1973",
    "A1977:H1995",
    "A1974:A1976",
    "A1997:H1999",
    "A1996",
    "C1996:H1996",
    "A2000",
    "C2000:H2000",
    "A2001:H2003",
    "A2004:A2007",
    "C1974:H1976",
    "C2004:H2006",
    "A2007:H2033",
    "C2034:H2035",
    "A2034:A2035",
    "A2036:H2040",
    "A2042:H2044",
    "A2041",
    "C2041:H2041",
    "A2046:H2049",
    "A2045",
    "C2045:H2045",
    "A2051:H2051",
    "A2050",
    "C2050:H2050",
    "A2053:H2056",
    "A2052",
    "C2052:H2052",
    "A2058:H2060",
    "A2057",
    "C2057:H2057",
    "A2061",
    "C2061:H2061",
    "A2062:H2066",
    "A2070:H2070",
    "A2067:A2069",
    "C2067:H2069",
    "A2073:H2073",
    "A2071:A2072",
    "C2071:H2072",
    "C2074:H2075",
    "A2074:A2075",
    "A2076:H2079",
    "A2080",
    "C2080:H2080",
    "A2081:H2085",
    "A2087:H2087",
    "A2086",
    "C2086:H2086",
    "A2089:H2090",
    "A2088",
    "C2088:H2088",
    "A2091:A2092",
    "C2091:H2092",
    "A1:H51",
    "A53:H1955",
    "A2093:H2097",
    "A2100:H2100",
    "A2098:A2099",
    "C2098:H2099",
    "A2101:A2102",
    "C2101:H2102",
    "A2103:H2111",
    "C2121:H2127",
    "A2124:B2124",
    "A2123",
    "A2126:B2127",
    "A2125",
    "A2121:B2122",
    "A2128:H2130",
    "A2113:H2116",
    "A2112",
    "C2112:H2112",
    "A2118:H2120",
    "A2117",
    "C2117:H2117",
    "A2139:H2150",
    "A2131:A2138",
    "C2131:H2138",
    "A2153:H2155",
    "A2151:A2152",
    "C2151:H2152",
    "A2156",
    "C2156:H2156",
    "A2157:H2160",
    "A2162:H2164",
    "A2161",
    "C2161:H2161",
    "C2165:H2167",
    "A2165:B2165",
    "A2167:B2167",
    "A2166",
    "A2168:H2170",
    "A2172:H2175",
    "A2171",
    "C2171:H2171",
    "A2177:H2180",
    "A2176",
    "C2176:H2176",
    "A2182:H2185",
    "A2181",
    "C2181:H2181",
    "A2186",
    "C2186:H2186",
    "A2187:H2187",
    "A2189:H2190",
    "A2188",
    "C2188:H2188",
    "A2191",
    "C2191:H2191",
    "A2200:H2202",
    "C2199:H2199",
    "C2203:H2205",
    "A2198:B2200",
    "A2192:H2198",
    "A2201:A2205",
    "A2204:B2210",
    "A2206:H2213",
    "A2215:H2218",
    "A2214",
    "C2214:H2214",
    "A2219",
    "C2219:H2219",
    "A2220:H2224",
    "A2225",
    "C2225:H2225",
    "C2229:H2229",
    "A2230:H2231",
    "A2232",
    "C2232:H2232",
    "A2233:H2235",
    "A2239:H2242",
    "C2236:H2238",
    "A2243:B2243",
    "A2244",
    "C2243:H2253",
    "A2254:H2256",
    "A2245:B2256",
    "A2236:A2239",
    "A2228:A2230",
    "A2226:H2228",
    "A2257:A2258",
    "C2257:H2258",
    "A2259:H2260",
    "A2262:H2262",
    "A2261",
    "C2261:H2261",
    "A2265:H2265",
    "C2263:H2264",
    "A2263:B2263",
    "A2264",
    "A2268:H2268",
    "C2266:H2267",
    "A2266:B2266",
    "A2267",
    "C2269:H2270",
    "A2269:B2269",
    "A2270",
    "A2271:H2271",
    "A2274:H2274",
    "A2272:A2273",
    "C2272:H2273",
    "C2275:H2287",
    "A2287:B2287",
    "A2285:B2285",
    "A2275:B2282",
    "A2283:A2286",
    "B2284",
    "A2288:H1048576"
)

$newApplies = $null
foreach ($addr in $cfRanges) {
    $r = $ws.Range($addr)
    if ($newApplies -eq $null) {
        $newApplies = $r
    } else {
        $newApplies = $excel.Union($newApplies, $r)
    }
}

$mainFC = $ws.Range("A1:H1048576").FormatConditions
for ($i = 1; $i -le $mainFC.Count; $i++) {
    $mainFC.Item($i).ModifyAppliesToRange($newApplies)
}

# --- Update the active selection to match where editing ended ---
$ws.Activate()
$ws.Range("B2290").Select()
try { $excel.ActiveWindow.ScrollRow = 2272 } catch {}
